$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose trajectory (columns K:AS, i.e. years 1-35) should be
# flattened to match the year-0 value already stored in column J.
$rows = @(3, 4, 5, 6, 9)

foreach ($r in $rows) {
    $baseValue = $ws.Cells.Item($r, 10).Value2  # column J = 10 (year 0)
    for ($c = 11; $c -le 45; $c++) {             # columns K(11) .. AS(45)
        $ws.Cells.Item($r, $c).Value2 = $baseValue
    }
}
